$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "60.027.84"
$ws.Range("D3").Value = "2.562.70"
$ws.Range("E3").Value = "  +8.70%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'510.78"
$ws.Range("E5").Value = "  +6.92%  "
$ws.Range("D6").Value = "'159.85"
$ws.Range("E6").Value = "  +8.41%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -3.10%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'0.990"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "2.612.21"
$ws.Range("E9").Value = "  +10.55%  "
$ws.Range("D10").Value = "'6.13"
$ws.Range("E10").Value = "  +12.23%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("D12").Value = "'0.345"
$ws.Range("E12").Value = "  +5.55%  "
$ws.Range("D13").Value = "'0.127"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "3.013.95"
$ws.Range("E14").Value = "  +9.11%  "
$ws.Range("D15").Value = "59.682.95"
$ws.Range("E15").Value = "  +7.99%  "
$ws.Range("D16").Value = "'22.07"
$ws.Range("E16").Value = "  +9.84%  "
$ws.Range("D17").Value = "'0.0000139"
$ws.Range("E17").Value = "  +6.99%  "
$ws.Range("D18").Value = "2.601.66"
$ws.Range("E18").Value = "  +10.49%  "
$ws.Range("D19").Value = "'4.80"
$ws.Range("E19").Value = "  +5.22%  "
$ws.Range("D20").Value = "'343.21"
$ws.Range("E20").Value = "  +8.52%  "
$ws.Range("D21").Value = "'10.48"
$ws.Range("E21").Value = "  +9.14%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  +8.33%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'60.14"
$ws.Range("E24").Value = "  +5.80%  "
$ws.Range("D25").Value = "'0.422"
$ws.Range("E25").Value = "  +6.60%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  +9.63%  "
$ws.Range("D27").Value = "2.686.34"
$ws.Range("E27").Value = "  +9.75%  "
$ws.Range("D28").Value = "'0.990"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "0.0₃0846"
$ws.Range("E29").Value = "  +12.68%  "
$ws.Range("D30").Value = "'7.41"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("D31").Value = "'0.996"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'157.67"
$ws.Range("E32").Value = "  +8.09%  "
$ws.Range("D33").Value = "'19.55"
$ws.Range("E33").Value = "  +7.45%  "
$ws.Range("D34").Value = "'1.58"
$ws.Range("E34").Value = "  +6.99%  "
$ws.Range("D35").Value = "'5.57"
$ws.Range("E35").Value = "  +8.73%  "
$ws.Range("D36").Value = "'1.21"
$ws.Range("E36").Value = "  +10.10%  "
$ws.Range("D37").Value = "'3.97"
$ws.Range("E37").Value = "  +10.64%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'314.43"
$ws.Range("E38").Value = "  +23.99%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.867"
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("D40").Value = "'3.77"
$ws.Range("E40").Value = "  +10.88%  "
$ws.Range("D41").Value = "'1.46"
$ws.Range("E41").Value = "  +9.21%  "
$ws.Range("D42").Value = "'35.30"
$ws.Range("E42").Value = "  +4.79%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.635"
$ws.Range("E43").Value = "  +9.99%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.102"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'0.0574"
$ws.Range("E45").Value = "  +10.63%  "
$ws.Range("D46").Value = "'0.788"
$ws.Range("E46").Value = "  +25.76%  "
$ws.Range("D47").Value = "'0.988"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").Value = "'5.02"
$ws.Range("E48").Value = "  +13.62%  "
$ws.Range("D49").Value = "'19.55"
$ws.Range("E49").Value = "  +16.61%  "
$ws.Range("D50").Value = "'0.0238"
$ws.Range("E50").Value = "  +7.52%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.000.35"
$ws.Range("E51").Value = "  +11.02%  "
